$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)
$full = $p1.Range
$insertStart = $full.End - 1

$r = $d.Range($insertStart, $insertStart)
$r.InsertAfter("https://github.com/darthrevan01/csd-310/tree/main/module-6")

$newRange = $d.Range($insertStart, $full.End - 1)
$newRange.Font.Name = "Times New Roman"
$newRange.Font.NameAscii = "Times New Roman"
$newRange.Font.NameBi = "Times New Roman"
$newRange.Font.Size = 12
$newRange.Font.SizeBi = 12
